$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) stays text so values like "1.000" or "0.00001046" are not
# reinterpreted/renormalized as numbers by Excel when assigned.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.322.35"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "1.775.50"
$ws.Range("E3").Value = "  +3.81%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "313.53"
$ws.Range("E5").Value = "  +2.20%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "0.5241"
$ws.Range("E7").Value = "  +10.89%  "
$ws.Range("D8").Value = "0.3620"
$ws.Range("E8").Value = "  +5.83%  "
$ws.Range("D9").Value = "42.53"
$ws.Range("E9").Value = "  +1.51%  "
$ws.Range("D10").Value = "0.07373"
$ws.Range("E10").Value = "  +1.61%  "
$ws.Range("D11").Value = "1.093"
$ws.Range("E11").Value = "  +5.00%  "
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("E13").Value = "  +3.89%  "
$ws.Range("D14").Value = "6.072"
$ws.Range("E14").Value = "  +3.95%  "
$ws.Range("D15").Value = "1.774.03"
$ws.Range("E15").Value = "  +4.16%  "
$ws.Range("D16").Value = "6.974"
$ws.Range("E16").Value = "  +2.23%  "
$ws.Range("D17").Value = "88.41"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").Value = "0.00001046"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("D19").Value = "0.06423"
$ws.Range("E19").Value = "  +1.09%  "
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "16.76"
$ws.Range("E21").Value = "  +1.90%  "
$ws.Range("D22").Value = "5.847"
$ws.Range("E22").Value = "  +4.63%  "
$ws.Range("D23").Value = "27.397.96"
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("D24").Value = "11.32"
$ws.Range("E24").Value = "  +4.59%  "
$ws.Range("D25").Value = "2.071"
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("D26").Value = "153.90"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("D27").Value = "20.15"
$ws.Range("E27").Value = "  +3.20%  "
$ws.Range("D28").Value = "2.356"
$ws.Range("E28").Value = "  +13.43%  "
$ws.Range("D29").Value = "1.976.06"
$ws.Range("E29").Value = "  +4.01%  "
$ws.Range("D30").Value = "121.27"
$ws.Range("E30").Value = "  +1.54%  "
$ws.Range("D31").Value = "1.063"
$ws.Range("E31").Value = "  +5.58%  "
$ws.Range("D32").Value = "0.09794"
$ws.Range("E32").Value = "  +7.13%  "
$ws.Range("D33").Value = "5.567"
$ws.Range("E33").Value = "  +5.12%  "
$ws.Range("D34").Value = "3.617"
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("E35").Value = "  +2.08%  "
$ws.Range("D36").Value = "0.05976"
$ws.Range("E36").Value = "  +2.89%  "
$ws.Range("D37").Value = "11.19"
$ws.Range("E37").Value = "  +1.37%  "
$ws.Range("D38").Value = "4.863"
$ws.Range("E38").Value = "  +2.82%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "0.6151"
$ws.Range("E39").Value = "  +5.17%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "0.2026"
$ws.Range("E40").Value = "  +1.93%  "
$ws.Range("D41").Value = "1.430"
$ws.Range("E41").Value = "  +2.88%  "
$ws.Range("D42").Value = "8.085"
$ws.Range("E42").Value = "  +8.63%  "
$ws.Range("D43").Value = "1.146"
$ws.Range("E43").Value = "  +4.16%  "
$ws.Range("D44").Value = "13.11"
$ws.Range("E44").Value = "  +3.66%  "
$ws.Range("D45").Value = "0.5773"
$ws.Range("E45").Value = "  +2.75%  "
$ws.Range("D46").Value = "3.630"
$ws.Range("E46").Value = "  +2.20%  "
$ws.Range("D47").Value = "121.30"
$ws.Range("D48").Value = "1.891"
$ws.Range("E48").Value = "  +3.09%  "
$ws.Range("D49").Value = "1.112"
$ws.Range("E49").Value = "  +2.75%  "
$ws.Range("D50").Value = "0.06715"
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("D51").Value = "70.72"
$ws.Range("E51").Value = "  +1.65%  "
